$wb = $excel.ActiveWorkbook

# ---- Seating sheet ----
$seating = $wb.Worksheets.Item("Seating")

# G8: "See Fu" -> "Chris Ng"
$seating.Range("G8").Value = "Chris Ng"

# Row 13: remove H13 (was "Chris Ng")
$seating.Range("H13").ClearContents()

# Row 18-20: add new F column entries
$seating.Range("F18").Value = "Bean Man"
$seating.Range("F19").Value = "Foo Kwai"
$seating.Range("F20").Value = "Ellen "

# G22: "Bean Man" -> "See Fu"
$seating.Range("G22").Value = "See Fu"

# Row 23-24: remove G23, G24 (moved to F19, F20)
$seating.Range("G23").ClearContents()
$seating.Range("G24").ClearContents()

# Update the active selection on Seating to F22
$seating.Activate()
$seating.Range("F22").Select()

# ---- Contact sheet ----
$contact = $wb.Worksheets.Item("Contact")

# B2: "C 吻" -> "Celia"
$contact.Range("B2").Value = "Celia"

# Update the active selection on Contact to E9
$contact.Activate()
$contact.Range("E9").Select()

# Restore Seating as the active/tabSelected sheet (matches original tabSelected="1")
$seating.Activate()
